$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sao Paulo
$ws.Range("B2").Value = 38.49536317209137
$ws.Range("D2").Value = 2024

# Row 3: Mato Grosso do Sul -> Sergipe
$ws.Range("A3").Value = "Sergipe"
$ws.Range("B3").Value = 22.62152262886292
$ws.Range("D3").Value = 2024

# Row 4: Mato Grosso
$ws.Range("B4").Value = 21.72911769654788
$ws.Range("D4").Value = 2024

# Row 5: Sergipe -> Mato Grosso do Sul
$ws.Range("A5").Value = "Mato Grosso do Sul"
$ws.Range("B5").Value = 20.03442551062107
$ws.Range("D5").Value = 2024

# Row 6: Rondonia
$ws.Range("B6").Value = 17.77559457850473
$ws.Range("D6").Value = 2024

# Row 7: Goias
$ws.Range("B7").Value = 17.31348793054583
$ws.Range("D7").Value = 2024

# Row 8: Brasil
$ws.Range("B8").Value = 16.66347210408774
$ws.Range("D8").Value = 2024

# Row 9: Nordeste
$ws.Range("B9").Value = 9.443802699703051
$ws.Range("D9").Value = 2024
